# Process-tree ("עץ תהליכים") deck — relabel four SmartArt boxes.
#
# Two pairs of node labels are swapped back to their "correct" wording:
#   "חיפוש עמדה"      <->  "רישום שחקן לעמדה"
#   "מחיקת טכנאי"     <->  "עדכון טכנאי"
#
# The SmartArt graphic keeps a rendered drawing cache (dsp:drawing) in
# sync with its data model, so editing a node's text through the
# SmartArt object model updates both ppt/diagrams/data1.xml and
# ppt/diagrams/drawing1.xml, matching the diff.

function Find-SmartArtShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasSmartArt) {
            return $sh
        }
        if ($sh.Type -eq 6) {
            # msoGroup -- the diagram lives inside a group shape on this slide.
            $found = Find-SmartArtShape $sh.GroupItems
            if ($found -ne $null) {
                return $found
            }
        }
    }
    return $null
}

function Find-SmartArtNodeByText($smartArt, [string]$text) {
    for ($i = 1; $i -le $smartArt.AllNodes.Count; $i++) {
        $node = $smartArt.AllNodes.Item($i)
        if ($node.TextFrame2.TextRange.Text -eq $text) {
            return $node
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$diagShape = Find-SmartArtShape $s.Shapes
$sa = $diagShape.SmartArt

# Resolve every node that needs to change BEFORE writing any new text,
# so the "current text" lookups stay unambiguous even though the two
# pairs are being swapped with each other.
$nodeSearchPos   = Find-SmartArtNodeByText $sa "חיפוש עמדה"
$nodeRegisterPos = Find-SmartArtNodeByText $sa "רישום שחקן לעמדה"
$nodeDeleteTech  = Find-SmartArtNodeByText $sa "מחיקת טכנאי"
$nodeUpdateTech  = Find-SmartArtNodeByText $sa "עדכון טכנאי"

$nodeSearchPos.TextFrame2.TextRange.Text   = "רישום שחקן לעמדה"
$nodeRegisterPos.TextFrame2.TextRange.Text = "חיפוש עמדה"
$nodeDeleteTech.TextFrame2.TextRange.Text  = "עדכון טכנאי"
$nodeUpdateTech.TextFrame2.TextRange.Text  = "מחיקת טכנאי"
